$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 68, shifting the existing rows 68-149 down to 69-150.
$ws.Rows("68:68").Insert("xlShiftDown")

# Populate the newly inserted row 68 with the new data record.
$ws.Range("A68").Value = 4
$ws.Range("B68").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C68").Value = "Los Lagos"
$ws.Range("D68").Value = [DateTime]"2021-12-09"
$ws.Range("E68").Value = 10
$ws.Range("F68").Value = 100112039
$ws.Range("G68").Value = "Ciboulette"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 120
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = 2500
$ws.Range("N68").Value = '$/docena de atados'
$ws.Range("O68").Value = "Región Metropolitana"
$ws.Range("P68").Value = 833
$ws.Range("Q68").Value = 3
$ws.Range("R68").Value = "Hortaliza"
